# Applies two content edits to the deck:
#  1. Slide 21 ("Ejercicio de ampliacion (IV)"): tweak wording of the second
#     bullet in the Content Placeholder from "...simplemente son implicitos"
#     to "...simplemente estan implicitos".
#  2. Slide 23 ("Depurando la aplicacion (II)"): append a new level-1 bullet
#     ("Y tambien podreis evaluar expresiones") after the last existing
#     bullet in the Content Placeholder.

$p = $ppt.ActivePresentation

# --- Edit 1: slide 21 wording tweak -----------------------------------
$slide21 = $p.Slides.Item(21)
$contentShape21 = $slide21.Shapes.Item(2)
$tr21 = $contentShape21.TextFrame.TextRange
$para21 = $tr21.Paragraphs(2)
$para21.Runs(1).Text = "Esto no quiere decir que los elementos que te acabamos de enseñar ya no existan – simplemente están implícitos"

# --- Edit 2: slide 23 new bullet ---------------------------------------
$slide23 = $p.Slides.Item(23)
$contentShape23 = $slide23.Shapes.Item(2)
$tr23 = $contentShape23.TextFrame.TextRange
$tr23.InsertAfter("`nY también podréis evaluar expresiones")
